$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated historical values and dividends for rows 2-18 (Weekly Performance %, Performance %, Value (£))
$updates = @(
    @{ Row = 2;  B = -10.05; C = -36.67; D = 657.7212071228028 },
    @{ Row = 3;  B = -0.89;  C = 11.07;  D = 1242.076779157063 },
    @{ Row = 4;  B = -4.85;  C = -40.42; D = 244.1154100140587 },
    @{ Row = 5;  B = -3.15;  C = -19.37; D = 207.1258462261812 },
    @{ Row = 6;  B = 3.86;   C = 127.79; D = 1284.492566079939 },
    @{ Row = 7;  B = 0.75;   C = -31.43; D = 498.751945613536 },
    @{ Row = 8;  B = -9.33;  C = -35.3;  D = 298.7874267198416 },
    @{ Row = 9;  B = -3.06;  C = -69.33; D = 370.2000045776367 },
    @{ Row = 10; B = 1.53;   C = -0.99;  D = 598.8 },
    @{ Row = 11; B = -1;     C = 4.51;   D = 533.0465272023926 },
    @{ Row = 12; B = -6.47;  C = 40.63;  D = 838.4581796676636 },
    @{ Row = 13; B = 3.94;   C = -13.49; D = 1289.729481161152 },
    @{ Row = 14; B = -6.03;  C = 30.87;  D = 229.7849179078067 },
    @{ Row = 15; B = -4.25;  C = -29.91; D = 488.4358902142849 },
    @{ Row = 16; B = -2.53;  C = 23.03;  D = 550.5192316595765 },
    @{ Row = 17; B = -8.1;   C = -73.67; D = 100.7946215176422 },
    @{ Row = 18; B = -1.52;  C = -5.67;  D = 9432.840034841578 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 2).Value = $u.B
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
}

# Row 19 only has C and D numeric (B19 is the "---" placeholder text)
$ws.Cells.Item(19, 3).Value = 1228.937934914572
$ws.Cells.Item(19, 4).Value = 1334.581746877552
